$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record (Blue Giant, dated 2022-03-07 / serial 44627) was
# reported ahead of the existing list, so insert a fresh row at position 152
# and push every row from the old 152 down to 200 (dimension A1:T199 -> A1:T200).
$ws.Rows.Item(152).Insert()

$ws.Range("A152").Value = 10
$ws.Range("B152").Value = 'Vega Modelo de Temuco'
$ws.Range("C152").Value = 'La Araucanía'
$ws.Range("D152").Value = 44627
$ws.Range("E152").Value = 9
$ws.Range("F152").Value = 'Fruta'
$ws.Range("G152").Value = 100103
$ws.Range("H152").Value = 'Frutos de hueso (carozo)'
$ws.Range("I152").Value = 100103002
$ws.Range("J152").Value = 'Ciruela'
$ws.Range("K152").Value = 'Blue Giant'
$ws.Range("L152").Value = 'Primera'
$ws.Range("M152").Value = 250
$ws.Range("N152").Value = 12000
$ws.Range("O152").Value = 12000
$ws.Range("P152").Value = 12000
$ws.Range("Q152").Value = '$/bandeja 18 kilos granel'
$ws.Range("R152").Value = "Región de O'Higgins"
$ws.Range("S152").Value = 667
$ws.Range("T152").Value = 18
